# Adds new ABP Object repository / Keyword Classes / ABP test case run
# timestamps by updating the "Date" column (column B) on the test-log
# sheets with the latest execution timestamps.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Mon Jun 02 21:54:47 IST 2025"
$ws.Range("B3").Value = "Mon Jun 02 21:55:43 IST 2025"
$ws.Range("B4").Value = "Mon Jun 02 21:56:25 IST 2025"
$ws.Range("B5").Value = "Mon Jun 02 21:57:07 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCCSCF")
$ws.Range("B2").Value = "Mon Jun 02 21:57:54 IST 2025"
$ws.Range("B3").Value = "Mon Jun 02 21:58:52 IST 2025"
$ws.Range("B4").Value = "Mon Jun 02 21:59:44 IST 2025"
$ws.Range("B5").Value = "Mon Jun 02 22:00:40 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCCDCF")
$ws.Range("B2").Value = "Mon Jun 02 22:01:33 IST 2025"
$ws.Range("B3").Value = "Mon Jun 02 22:02:36 IST 2025"
$ws.Range("B4").Value = "Mon Jun 02 22:03:32 IST 2025"
$ws.Range("B5").Value = "Mon Jun 02 22:04:24 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyAmount")
$ws.Range("B2").Value = "Mon Jun 02 22:06:38 IST 2025"

$ws = $wb.Worksheets.Item("OverUnderPay")
$ws.Range("B2").Value = "Mon Jun 02 22:10:59 IST 2025"
$ws.Range("B3").Value = "Mon Jun 02 22:11:34 IST 2025"

$ws = $wb.Worksheets.Item("NoOverPay")
$ws.Range("B2").Value = "Mon Jun 02 22:15:07 IST 2025"
